$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.700.57"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.43"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.44"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.15"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.02"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.975"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.38"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.656.66"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.309.07"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.649.39"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.19"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.51"
$ws.Range("E23").Value = "  +9.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  +21.43%  "
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.95"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.18"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.15"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0878"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +6.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.59"
$ws.Range("E36").Value = "  -10.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0378"
$ws.Range("E37").Value = "  +7.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.66"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("E41").Value = "  +4.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.68"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.05"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.01"
$ws.Range("E47").Value = "  +8.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.27"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.30"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.604.46"
$ws.Range("E51").Value = "  +4.63%  "
